$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.687.04'
$ws.Range('E2').Value = '  -0.97%  '

# Row 3
$ws.Range('D3').Value = '1.630.09'
$ws.Range('E3').Value = '  -1.01%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.00'
$ws.Range('E5').Value = '  -1.23%  '

# Row 6
$ws.Range('E6').Value = '  -0.99%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.13%  '

# Row 8
$ws.Range('E8').Value = '  -0.91%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0637'
$ws.Range('E9').Value = '  -1.19%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.45'
$ws.Range('E10').Value = '  -5.94%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  +0.09%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.641.11'
$ws.Range('E12').Value = '  -0.32%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.24'
$ws.Range('E13').Value = '  -1.04%  '

# Row 14
$ws.Range('D14').Value = '1.854.67'
$ws.Range('E14').Value = '  -0.96%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.551'
$ws.Range('E15').Value = '  -2.07%  '

# Row 16
$ws.Range('D16').Value = '0.0₃0766'
$ws.Range('E16').Value = '  -0.92%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.10'
$ws.Range('E17').Value = '  -0.54%  '

# Row 18
$ws.Range('D18').Value = '25.714.54'
$ws.Range('E18').Value = '  -0.86%  '

# Row 19
$ws.Range('E19').Value = '  -0.09%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  +0.94%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '194.01'
$ws.Range('E21').Value = '  +0.02%  '

# Row 22
$ws.Range('E22').Value = '  -0.49%  '

# Row 23
$ws.Range('E23').Value = '  +0.99%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.08%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.79'
$ws.Range('E25').Value = '  -0.98%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.32'
$ws.Range('E26').Value = '  -0.64%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.120'
$ws.Range('E27').Value = '  -3.55%  '

# Row 28
$ws.Range('E28').Value = '  -0.47%  '

# Row 29
$ws.Range('E29').Value = '  -0.62%  '

# Row 30
$ws.Range('E30').Value = '  -1.28%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0484'
$ws.Range('E31').Value = '  -2.74%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  +0.65%  '

# Row 33
$ws.Range('E33').Value = '  -0.11%  '

# Row 34
$ws.Range('E34').Value = '  +0.19%  '

# Row 35
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.894'
$ws.Range('E36').Value = '  -1.49%  '

# Row 37
$ws.Range('E37').Value = '  -0.79%  '

# Row 38
$ws.Range('E38').Value = '  -2.67%  '

# Row 39
$ws.Range('D39').Value = '1.103.22'
$ws.Range('E39').Value = '  -2.75%  '

# Row 40
$ws.Range('E40').Value = '  -1.01%  '

# Row 41
$ws.Range('E41').Value = '  +0.33%  '

# Row 42
$ws.Range('E42').Value = '  +0.81%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.78'
$ws.Range('E43').Value = '  +1.00%  '

# Row 44
$ws.Range('E44').Value = '  -1.28%  '

# Row 45
$ws.Range('D45').Value = '1.762.85'
$ws.Range('E45').Value = '  -0.95%  '

# Row 46
$ws.Range('D46').Value = '0.0₆0107'
$ws.Range('E46').Value = '  -2.54%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.96'
$ws.Range('E47').Value = '  -1.60%  '

# Row 48
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.418'
$ws.Range('E48').Value = '  -2.44%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.68'
$ws.Range('E49').Value = '  -1.14%  '

# Row 50
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.37'
$ws.Range('E50').Value = '  +3.94%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0502'
$ws.Range('E51').Value = '  -0.49%  '
